$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 45-46 (existing rows 45..56 shift down to 47..58)
$ws.Range("A45:A46").EntireRow.Insert()

# New row 45: weekly Sandia "Primera" price data for Peru
$ws.Cells.Item(45, 1).Value = 8
$ws.Cells.Item(45, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(45, 3).Value = "Coquimbo"
$ws.Cells.Item(45, 4).Value = 44522
$ws.Cells.Item(45, 5).Value = 4
$ws.Cells.Item(45, 6).Value = 100112028
$ws.Cells.Item(45, 7).Value = "Sandia"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 600
$ws.Cells.Item(45, 11).Value = 900
$ws.Cells.Item(45, 12).Value = 1000
$ws.Cells.Item(45, 13).Value = 950
$ws.Cells.Item(45, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(45, 15).Value = "Perú"
$ws.Cells.Item(45, 16).Value = 950
$ws.Cells.Item(45, 17).Value = 1
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# New row 46: weekly Sandia "Segunda" price data for Peru
$ws.Cells.Item(46, 1).Value = 8
$ws.Cells.Item(46, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(46, 3).Value = "Coquimbo"
$ws.Cells.Item(46, 4).Value = 44522
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(46, 6).Value = 100112028
$ws.Cells.Item(46, 7).Value = "Sandia"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Segunda"
$ws.Cells.Item(46, 10).Value = 400
$ws.Cells.Item(46, 11).Value = 700
$ws.Cells.Item(46, 12).Value = 800
$ws.Cells.Item(46, 13).Value = 750
$ws.Cells.Item(46, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(46, 15).Value = "Perú"
$ws.Cells.Item(46, 16).Value = 750
$ws.Cells.Item(46, 17).Value = 1
$ws.Cells.Item(46, 18).Value = "Hortaliza"
